$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 6-9: switch title (column A) cells to the
# wrapped-text style used throughout the rest of column A ---
$ws.Range("A6").WrapText = $true
$ws.Range("A7").WrapText = $true
$ws.Range("A8").WrapText = $true
$ws.Range("A9").WrapText = $true

# Row 9 grows to a two-line title, so its row height increases to 30
# (matching the other multi-line rows).
$ws.Rows.Item(9).RowHeight = 30

# --- Add the new paper as row 10 ---
$ws.Range("A10").Value = "The Importance of Temporal Resolution in Evaluating Residential Energy Storage"
$ws.Range("B10").Value = "https://ieeexplore.ieee.org/stamp/stamp.jsp?tp=&arnumber=8274019&tag=1"

# Hook up the hyperlink for the new row's link cell.
$ws.Hyperlinks.Add($ws.Range("B10"), "https://ieeexplore.ieee.org/stamp/stamp.jsp?tp=&arnumber=8274019&tag=1")

# Re-apply the same visual formatting (wrap text) used by the other rows;
# Hyperlinks.Add resets formatting on the target cell, so this must come
# after the call above.
$ws.Range("A10").WrapText = $true
$ws.Range("B10").WrapText = $true

# New row mirrors the other multi-line rows' height.
$ws.Rows.Item(10).RowHeight = 30
